$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Ke%Hu%NULL%0, Wei-jie%Guan%NULL%0, Ying%Bi%NULL%1, Wei%Zhang%NULL%0, Lanjuan%Li%NULL%0, Boli%Zhang%NULL%1, Qingquan%Liu%NULL%0, Yuanlin%Song%NULL%1, Xingwang%Li%NULL%0, Zhongping%Duan%NULL%1, Qingshan%Zheng%NULL%1, Zifeng%Yang%NULL%1, Jingyi%Liang%NULL%1, Mingfeng%Han%NULL%0, Lianguo%Ruan%NULL%2, Chaomin%Wu%NULL%1, Yunting%Zhang%NULL%1, Zhen-hua%Jia%NULL%1, Nan-shan%Zhong%NULL%0]'
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = 'Elsevier GmbH.'
